$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text / link / non-ambiguous numeric-looking price updates
$ws.Range("D2").Value = "54.224.91"
$ws.Range("E2").Value = "  -3.52%  "
$ws.Range("D3").Value = "2.271.21"
$ws.Range("E3").Value = "  -4.05%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("E5").Value = "  -2.76%  "
$ws.Range("E6").Value = "  -2.26%  "
$ws.Range("E7").Value = "  +0.16%  "
$ws.Range("E8").Value = "  -2.62%  "
$ws.Range("D9").Value = "2.268.59"
$ws.Range("E9").Value = "  -4.42%  "
$ws.Range("E10").Value = "  -5.42%  "
$ws.Range("E11").Value = "  +0.07%  "
$ws.Range("E12").Value = "  -0.39%  "
$ws.Range("E13").Value = "  -4.22%  "
$ws.Range("D14").Value = "2.670.25"
$ws.Range("E14").Value = "  -4.16%  "
$ws.Range("E15").Value = "  -1.22%  "
$ws.Range("D16").Value = "54.140.88"
$ws.Range("E16").Value = "  -3.65%  "
$ws.Range("E17").Value = "  -3.31%  "
$ws.Range("D18").Value = "2.292.91"
$ws.Range("E18").Value = "  -2.35%  "
$ws.Range("E19").Value = "  -2.02%  "
$ws.Range("E20").Value = "  +0.08%  "
$ws.Range("E21").Value = "  -3.79%  "
$ws.Range("E22").Value = "  -0.18%  "
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("E24").Value = "  -2.67%  "
$ws.Range("E25").Value = "  +0.33%  "
$ws.Range("E26").Value = "  +0.42%  "
$ws.Range("B27").Value = "WrappedeETH"
$ws.Range("C27").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D27").Value = "2.370.32"
$ws.Range("E27").Value = "  -4.46%  "
$ws.Range("B28").Value = "Kaspa"
$ws.Range("C28").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("E28").Value = "  -0.10%  "
$ws.Range("E29").Value = "  -1.03%  "
$ws.Range("E30").Value = "  -5.83%  "
$ws.Range("E31").Value = "  -2.95%  "
$ws.Range("D32").Value = "0.0₃0680"
$ws.Range("E32").Value = "  -4.47%  "
$ws.Range("E33").Value = "  -0.03%  "
$ws.Range("E34").Value = "  -0.73%  "
$ws.Range("E35").Value = "  +0.13%  "
$ws.Range("E37").Value = "  -0.61%  "
$ws.Range("E38").Value = "  +0.31%  "
$ws.Range("E39").Value = "  +1.34%  "
$ws.Range("E40").Value = "  -1.34%  "
$ws.Range("E41").Value = "  -2.22%  "
$ws.Range("E42").Value = "  +0.77%  "
$ws.Range("E43").Value = "  +0.26%  "
$ws.Range("E44").Value = "  -1.21%  "
$ws.Range("E45").Value = "  +0.63%  "
$ws.Range("E46").Value = "  +1.60%  "
$ws.Range("E47").Value = "  -1.00%  "
$ws.Range("E48").Value = "  +2.30%  "
$ws.Range("E49").Value = "  -3.32%  "
$ws.Range("E50").Value = "  -0.74%  "
$ws.Range("E51").Value = "  -1.63%  "

# Price cells whose new text would otherwise be auto-coerced to a Number by Excel
# (single "." decimal look-alikes) -- force Text type, then restore the default
# "Normal" style so no stray style index is left on the cell.
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "491.73"
$c.Style = "Normal"
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "127.15"
$c.Style = "Normal"
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.529"
$c.Style = "Normal"
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.323"
$c.Style = "Normal"
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "21.48"
$c.Style = "Normal"
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "9.80"
$c.Style = "Normal"
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "297.88"
$c.Style = "Normal"
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "6.25"
$c.Style = "Normal"
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "63.80"
$c.Style = "Normal"
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "0.147"
$c.Style = "Normal"
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "7.11"
$c.Style = "Normal"
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "162.99"
$c.Style = "Normal"
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "5.80"
$c.Style = "Normal"
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "0.998"
$c.Style = "Normal"
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "17.43"
$c.Style = "Normal"
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "1.19"
$c.Style = "Normal"
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "3.63"
$c.Style = "Normal"
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "35.47"
$c.Style = "Normal"
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "126.14"
$c.Style = "Normal"
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "242.18"
$c.Style = "Normal"
